# updated test cases for share skill and manage listings
$wb = $excel.ActiveWorkbook

# ---- ShareSkill: add a new test-data row (row 3) ----
$wsShare = $wb.Worksheets.Item("ShareSkill")

# Clone row 2's formatting onto row 3 so styles line up exactly.
$wsShare.Range("A2:N2").Copy()
$wsShare.Range("A3:N3").PasteSpecial(-4122)
# K3 should use the hh:mm:ss style (same as column J), not hh:mm like K2.
$wsShare.Range("J2").Copy()
$wsShare.Range("K3").PasteSpecial(-4122)

$wsShare.Cells.Item(3, 1).Value = "Selenium"
$wsShare.Cells.Item(3, 2).Value = "Selenium"
$wsShare.Cells.Item(3, 3).Value = "Business"
$wsShare.Cells.Item(3, 4).Value = "Other"
$wsShare.Cells.Item(3, 5).Value = "Test"
$wsShare.Cells.Item(3, 6).Value = "Hourly basis service"
$wsShare.Cells.Item(3, 7).Value = "Online"
$wsShare.Cells.Item(3, 8).Value = 44585
$wsShare.Cells.Item(3, 9).Value = 44590
$wsShare.Cells.Item(3, 10).Value = 0.6326388888888889
$wsShare.Cells.Item(3, 11).Value = 0.6326388888888889
$wsShare.Cells.Item(3, 12).Value = "Credit"
$wsShare.Cells.Item(3, 13).Value = 5
$wsShare.Cells.Item(3, 14).Value = "Active"

# ---- ManageListings: add a new test-data row (row 3) ----
$wsManage = $wb.Worksheets.Item("ManageListings")

$wsManage.Range("A2:N2").Copy()
$wsManage.Range("A3:N3").PasteSpecial(-4122)
$wsManage.Range("J2").Copy()
$wsManage.Range("K3").PasteSpecial(-4122)

$wsManage.Cells.Item(3, 1).Value = "Manage"
$wsManage.Cells.Item(3, 2).Value = "ManageListings"
$wsManage.Cells.Item(3, 3).Value = "Digital Marketing"
$wsManage.Cells.Item(3, 4).Value = "Video Marketing"
$wsManage.Cells.Item(3, 5).Value = "Test"
$wsManage.Cells.Item(3, 6).Value = "Hourly basis service"
$wsManage.Cells.Item(3, 7).Value = "Online"
$wsManage.Cells.Item(3, 8).Value = 44585
$wsManage.Cells.Item(3, 9).Value = 44590
$wsManage.Cells.Item(3, 10).Value = 0.6326388888888889
$wsManage.Cells.Item(3, 11).Value = 0.6326388888888889
$wsManage.Cells.Item(3, 12).Value = "Credit"
$wsManage.Cells.Item(3, 13).Value = 5
$wsManage.Cells.Item(3, 14).Value = "Active"
